# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.264.55"
$ws.Range("E2").Value = "  -0.52%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.413.21"
$ws.Range("E3").Value = "  +1.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.89"
$ws.Range("E5").Value = "  -1.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "661.76"
$ws.Range("E6").Value = "  +1.67%  "

# Row 7
$ws.Range("E7").Value = "  -6.05%  "

# Row 8
$ws.Range("E8").Value = "  -6.24%  "

# Row 9
$ws.Range("E9").Value = "  -3.22%  "

# Row 10
$ws.Range("E10").Value = "  +0.04%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.412.70"
$ws.Range("E11").Value = "  +1.91%  "

# Row 12
$ws.Range("E12").Value = "  +2.92%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.63"
$ws.Range("E13").Value = "  -2.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.56"
$ws.Range("E14").Value = "  +17.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "97.989.59"
$ws.Range("E15").Value = "  -1.37%  "

# Row 16
$ws.Range("E16").Value = "  -0.92%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.041.15"
$ws.Range("E17").Value = "  +1.46%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.00"
$ws.Range("E18").Value = "  +20.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.413.86"
$ws.Range("E19").Value = "  +1.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.568"
$ws.Range("E20").Value = "  +30.40%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.69"
$ws.Range("E21").Value = "  +4.40%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.05"
$ws.Range("E22").Value = "  +7.21%  "

# Row 23
$ws.Range("E23").Value = "  -3.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "510.38"
$ws.Range("E24").Value = "  -5.20%  "

# Row 25
$ws.Range("E25").Value = "  -3.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.75"
$ws.Range("E26").Value = "  +7.66%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "100.64"
$ws.Range("E27").Value = "  -2.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.91"
$ws.Range("E28").Value = "  +1.75%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.598.22"
$ws.Range("E29").Value = "  +2.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.150"
$ws.Range("E30").Value = "  -1.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.61"
$ws.Range("E31").Value = "  +6.28%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.199"
$ws.Range("E32").Value = "  +5.24%  "

# Row 33
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.02%  "

# Row 35
$ws.Range("B35").Value = "PolygonEcosystemToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.576"
$ws.Range("E35").Value = "  +6.95%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  +14.57%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.92"
$ws.Range("E37").Value = "  +2.40%  "

# Row 38
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.52"
$ws.Range("E38").Value = "  +14.71%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.99"
$ws.Range("E39").Value = "  +2.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "538.88"
$ws.Range("E40").Value = "  +4.12%  "

# Row 41
$ws.Range("E41").Value = "  -1.91%  "

# Row 42
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.880"
$ws.Range("E43").Value = "  +7.13%  "

# Row 44
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.16"
$ws.Range("E44").Value = "  +16.88%  "

# Row 45
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.71"
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0437"
$ws.Range("E46").Value = "  +8.19%  "

# Row 47
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.86"
$ws.Range("E47").Value = "  +16.40%  "

# Row 48
$ws.Range("E48").Value = "  +15.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").Value = "  -2.08%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.30"
$ws.Range("E50").Value = "  -2.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.86"
$ws.Range("E51").Value = "  +8.57%  "
